# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# The workbook currently ends with a single "总计" (grand-total) summary
# sheet (tab #6). We:
#   1. Turn that sheet into the new "2022-Q1" per-fund holdings sheet
#      (same shape/columns as the other quarterly sheets).
#   2. Append a brand-new "总计" sheet right after it, holding the same
#      summary table as before plus a new leading row for 2022-Q1.
# ---------------------------------------------------------------------------

function Set-TextValue($ws, $stage, $targetAddr, $text) {
    # Writing a numeric-looking string straight into a General-formatted
    # cell (e.g. Range.Value = "001167") gets auto-coerced to a number and
    # loses leading zeros / introduces float noise (e.g. "4.03" -> 4.0300000000000002).
    # Stage the text in a scratch cell formatted as Text, then copy only the
    # *value* over so the destination cell keeps its original (unstyled)
    # formatting while still storing a genuine text value.
    $stage.NumberFormat = "@"
    $stage.Value = $text
    $stage.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

$wsQ1 = $wb.Worksheets.Item(6)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------
# 1) Rebuild "2022-Q1" fund-holdings sheet
# ---------------------------------------------------------------------

$stage1 = $wsQ1.Range("Z100")

# Header row (B1:D1 already carry the bold/bordered header style; extend
# that same style across the new E1:H1 header cells)
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

$wsQ1.Range("D1").Copy()
$wsQ1.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fund rows: code, name, fund size, stock position, position %,
# holding value, position rank
$q1Data = @(
    @("001167", "金鹰科技创新股票",   "4.03", "94.55", "5.14", "0.2071", 3),
    @("210009", "金鹰核心资源混合",   "3.86", "94.96", "4.78", "0.1845", 4),
    @("162102", "金鹰中小盘精选混合", "4.60", "76.52", "3.70", "0.1702", 6),
    @("000458", "英大领先回报混合",   "1.11", "93.55", "3.11", "0.0345", 1),
    @("001270", "英大灵活配置混合A", "0.56", "93.18", "3.06", "0.0171", 1),
    @("001271", "英大灵活配置混合B", "0.28", "93.18", "3.06", "0.0086", 1)
)

$r = 2
foreach ($row in $q1Data) {
    $wsQ1.Range("A$r").Value = ($r - 2)
    Set-TextValue $wsQ1 $stage1 "B$r" $row[0]
    $wsQ1.Range("C$r").Value = $row[1]
    Set-TextValue $wsQ1 $stage1 "D$r" $row[2]
    Set-TextValue $wsQ1 $stage1 "E$r" $row[3]
    Set-TextValue $wsQ1 $stage1 "F$r" $row[4]
    Set-TextValue $wsQ1 $stage1 "G$r" $row[5]
    $wsQ1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

$stage1.Clear()
$excel.CutCopyMode = $false

# Extend the column-A index style (bold/bordered/centered) down through
# the new row 7
$wsQ1.Range("A2").Copy()
$wsQ1.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Build the new "总计" summary sheet
# ---------------------------------------------------------------------

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 6, 0.62),
    @("2021-Q4", 21, 7.92),
    @("2021-Q3", 1, 0.02),
    @("2021-Q2", 37, 13.86),
    @("2021-Q1", 11, 0.64),
    @("2020-Q4", 4, 0.05)
)

$r = 2
foreach ($row in $totalData) {
    $wsTotal.Range("A$r").Value = ($r - 2)
    $wsTotal.Range("B$r").Value = $row[0]
    $wsTotal.Range("C$r").Value = $row[1]
    $wsTotal.Range("D$r").Value = $row[2]
    $r = $r + 1
}

# Match header (B1:D1) and column-A index styling against the sheets that
# already carry the correct bold/bordered/centered look
$wsQ1.Range("D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ1.Range("A2").Copy()
$wsTotal.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
